# "duplicate match issue fixed"
# The match sheet had been populated with data from a duplicate/incorrect
# match. This corrects the Bangladesh batting card (A:F), the Mumbai
# India batting card (J:O), the innings totals (row 16) and both bowling
# figures tables (rows 21-25) to the real scorecard values.
#
# Note: a handful of cells in the "Overs" columns (B/K on rows 21-25 and
# C16/L16) hold cricket-over notation ("3.0", "1.3", ...) stored as TEXT
# in the workbook, not numbers (so "1.3" isn't coerced/rounded as 1.3
# decimal). Assigning a bare numeric-looking string to `.Value` would get
# auto-converted back to a number by Excel, so those are written with a
# leading apostrophe to force text, matching the original representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bangladesh batting (A:F) / Mumbai India batting (J:O) — row 2 ---
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "LBW"
$ws.Range("K2").Value = 39
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = "Bowled"
$ws.Range("N2").Value = " Taskin Ahmed"

# --- row 3 ---
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "Caught"
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = "LBW"
$ws.Range("N3").Value = " Mustafizur Rahman"

# --- row 4 ---
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("L4").Value = 4
$ws.Range("N4").Value = " Mahedi Hasan"

# --- row 5 ---
$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 8
$ws.Range("E5").Value = " Hardik Pandya"
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 3
$ws.Range("N5").Value = " Shamim Hossain"

# --- row 6 ---
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 16
$ws.Range("E6").Value = " Jasprit Bumrah"
$ws.Range("K6").Value = 15
$ws.Range("L6").Value = 6
$ws.Range("N6").Value = " Mustafizur Rahman"

# --- row 7 ---
$ws.Range("B7").Value = 108
$ws.Range("C7").Value = 33
$ws.Range("D7").Value = "NOT OUT"
$ws.Range("E7").Value = " "
$ws.Range("K7").Value = 30
$ws.Range("L7").Value = 9
$ws.Range("M7").Value = "Caught"
$ws.Range("N7").Value = " Shoriful Islam"

# --- row 8 ---
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "LBW"
$ws.Range("E8").Value = " Jasprit Bumrah"
$ws.Range("K8").Value = 28
$ws.Range("L8").Value = 9
$ws.Range("N8").Value = " Mahedi Hasan"

# --- row 9 ---
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = "Caught"
$ws.Range("E9").Value = " Bhuvneshwar Kumar"
$ws.Range("K9").Value = 46
$ws.Range("L9").Value = 14
$ws.Range("M9").Value = "Bowled"

# --- row 10 ---
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = "Bowled"
$ws.Range("E10").Value = " Yuzvendra Chahal"
$ws.Range("K10").Value = 11
$ws.Range("L10").Value = 7
$ws.Range("M10").Value = "Caught"
$ws.Range("N10").Value = " Mustafizur Rahman"

# --- row 11 ---
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 3
$ws.Range("E11").Value = " Jasprit Bumrah"
$ws.Range("K11").Value = 25
$ws.Range("L11").Value = 8
$ws.Range("M11").Value = "NOT OUT"
$ws.Range("N11").Value = " "

# --- row 12 ---
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = "Bowled"
$ws.Range("E12").Value = " Jasprit Bumrah"
$ws.Range("K12").Value = 10
$ws.Range("L12").Value = 6
$ws.Range("M12").Value = "Caught"
$ws.Range("N12").Value = " Mahedi Hasan"

# --- innings totals — row 16 ---
$ws.Range("A16").Value = 279
$ws.Range("C16").Value = "'16.3"
$ws.Range("D16").Value = 99
$ws.Range("J16").Value = 221
$ws.Range("L16").Value = "'13.3"
$ws.Range("M16").Value = 81

# --- bowling figures — row 21 ---
$ws.Range("A21").Value = "Mohommad Shami"
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 14.33
$ws.Range("K21").Value = "'2.0"
$ws.Range("L21").Value = 28
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 14

# --- row 22 ---
$ws.Range("A22").Value = "Bhuvneshwar Kumar"
$ws.Range("C22").Value = 59
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 19.67
$ws.Range("K22").Value = "'3.0"
$ws.Range("L22").Value = 65
$ws.Range("M22").Value = 2
$ws.Range("N22").Value = 21.67

# --- row 23 ---
$ws.Range("A23").Value = "Hardik Pandya"
$ws.Range("C23").Value = 50
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 16.67
$ws.Range("K23").Value = "'3.0"
$ws.Range("L23").Value = 44
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 14.67

# --- row 24 ---
$ws.Range("A24").Value = "Yuzvendra Chahal"
$ws.Range("B24").Value = "'4.0"
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 15.75
$ws.Range("K24").Value = "'3.0"
$ws.Range("L24").Value = 52
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 17.33

# --- row 25 ---
$ws.Range("A25").Value = "Jasprit Bumrah"
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 19.39
$ws.Range("K25").Value = "'2.3"
$ws.Range("L25").Value = 32
$ws.Range("N25").Value = 13.91
